# next/previous navigation added (testing)
# Update the "next" message-type row to describe the new "newFollow" message,
# adding a "direction (forward/backward)" data column ahead of the existing
# "previously following" / "now following" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 1).Value = "newFollow"
$ws.Cells.Item(10, 2).Value = "direction (forward/backward)"
$ws.Cells.Item(10, 3).Value = "previously following"
$ws.Cells.Item(10, 4).Value = "now following"

# Column B needs to widen to fit the new, longer text (mirrors Excel's
# "bestFit" autosizing behavior for this column, target stored width 24.15625).
$ws.Columns.Item(2).ColumnWidth = 23.33

$ws.Range("D10").Select()
